$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 20.29529466666667
$ws.Range("N2").Value = 60.885884
$ws.Range("O2").Value = 0.4032332285476398
$ws.Range("P2").Value = 0.4032332285476398
$ws.Range("Q2").Value = 7.213773066516445
$ws.Range("R2").Value = 64.923957598648
$ws.Range("S2").Value = 0.4032332285476398
$ws.Range("T2").Value = 0.4032332285476398

# Row 3
$ws.Range("O3").Value = 0.1953894087318433
$ws.Range("P3").Value = 0.1953894087318433
$ws.Range("S3").Value = 0.1953894087318433
$ws.Range("T3").Value = 0.1953894087318433

# Row 4
$ws.Range("M4").Value = 11.81535133333333
$ws.Range("N4").Value = 35.446054
$ws.Range("O4").Value = 0.2347510761885954
$ws.Range("P4").Value = 0.2347510761885954
$ws.Range("Q4").Value = 4.199656354820888
$ws.Range("R4").Value = 37.79690719338799
$ws.Range("S4").Value = 0.2347510761885954
$ws.Range("T4").Value = 0.2347510761885954

# Row 5
$ws.Range("M5").Value = 8.386535
$ws.Range("N5").Value = 25.159605
$ws.Range("O5").Value = 0.1666262865319216
$ws.Range("P5").Value = 0.1666262865319216
$ws.Range("Q5").Value = 2.980915591423333
$ws.Range("R5").Value = 26.82824032281
$ws.Range("S5").Value = 0.1666262865319216
$ws.Range("T5").Value = 0.1666262865319216
